$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-7 (columns D, K, L, M, N, O, P, R, S) per the diff.
# This is a row-wise permutation of the original data block.

$ws.Range("D2").Value2 = 44355
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1139

$ws.Range("D3").Value2 = 44342
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1361

$ws.Range("D4").Value2 = 44313
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1194

$ws.Range("D5").Value2 = 44699
$ws.Range("K5").Value = "Mankaki"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 29000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1639

$ws.Range("D6").Value2 = 44301
$ws.Range("K6").Value = "Hachiya"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1139

$ws.Range("D7").Value2 = 44305
$ws.Range("K7").Value = "Mankaki"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1361
